$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4: title/link update
$ws.Range("D4").Value = "#07-Pandas(판다스) Groupby와 Pivot table"
$ws.Range("E4").Value = "https://teddylee777.github.io/pandas/pandas-tutorial-07"

# Row 9: title/link update
$ws.Range("D9").Value = "PDSI – Edu 면접 가이드"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/pdsi-edu-interview-guide/#utm_source=rss&utm_medium=rss&utm_campaign=pdsi-edu-interview-guide"

# Row 29: title/link update
$ws.Range("D29").Value = "도커 컨테이너를 실행한 사용자 식별 방법 연구"
$ws.Range("E29").Value = "https://blog.promedius.ai/identify-docker-container-runner/"
